# Apply edits described by the commit:
#  "Added Password forgot option, added manage connection option,
#   added the UI for manage connection file"
#
# Concretely, on the "demo_file" worksheet this updates the DB-connection
# template text in column C (adds sourceuser / targetuser lines, fixes
# "targetdbType" -> "TargetdbType") and unifies column D's
# Inventory placeholder to "src_inventory:dest_inventory" for every data
# row, then adjusts row heights / the active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

$newDbDetails = "sourcedbType:mysql;`nsourceServer:localhost;`nsourcedb:source_db;`nsourceuser:acciom_user;`nTargetdbType:mysql;`ntargetdb:dest_db;`ntargetServer:localhost;`ntargetuser:Acciom_user;"
$newInventory = "src_inventory:dest_inventory"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $newDbDetails
    $ws.Cells.Item($r, 4).Value2 = $newInventory
}

# Row heights (points) - reflects the extra wrapped lines now in column C
$ws.Rows.Item(2).RowHeight = 142
$ws.Rows.Item(3).RowHeight = 126.1
$ws.Rows.Item(4).RowHeight = 126.1
$ws.Rows.Item(5).RowHeight = 136.35
$ws.Rows.Item(6).RowHeight = 114.9
$ws.Rows.Item(7).RowHeight = 114.9
$ws.Rows.Item(8).RowHeight = 114.9
$ws.Rows.Item(9).RowHeight = 146.55
$ws.Rows.Item(10).RowHeight = 136.35

# Update the saved view state: scrolled so A9 is the top-left visible cell,
# with C19 as the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("C19").Select() | Out-Null
